# Update cryptos list values (Price column D, Volume(1h) column E)
# Commit: Updated cryptos list on Thu Oct 17 13:55:54 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "67.187.45"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.609.56"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.98"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.77"
$ws.Range("E6").Value = "  -0.98%  "

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8: XRP
$ws.Range("E8").Value = "  +1.07%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.606.07"
$ws.Range("E9").Value = "  +0.08%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -3.18%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +0.38%  "

# Row 12: Toncoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.19"
$ws.Range("E12").Value = "  -1.32%  "

# Row 13: Cardano
$ws.Range("E13").Value = "  -2.76%  "

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.64"
$ws.Range("E14").Value = "  -0.50%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.087.02"
$ws.Range("E15").Value = "  -0.05%  "

# Row 16: ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  -3.86%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "66.960.49"
$ws.Range("E17").Value = "  -0.85%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.614.11"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "367.26"
$ws.Range("E19").Value = "  +0.44%  "

# Row 20: Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.05"
$ws.Range("E20").Value = "  -2.15%  "

# Row 21: Uniswap
$ws.Range("E21").Value = "  -4.42%  "

# Row 22: Polkadot
$ws.Range("E22").Value = "  -0.31%  "

# Row 23: SuiNetwork
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  +1.91%  "

# Row 24: Dai
$ws.Range("E24").Value = "  +0.33%  "

# Row 25: Aptos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.11"
$ws.Range("E25").Value = "  +1.78%  "

# Row 26: Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.80"
$ws.Range("E26").Value = "  -1.57%  "

# Row 27: WrappedeETH
$ws.Range("E27").Value = "  -0.21%  "

# Row 28: Bittensor
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "586.28"
$ws.Range("E28").Value = "  +2.00%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.48%  "

# Row 30: PEPE
$ws.Range("E30").Value = "  -3.09%  "

# Row 31: Fetch.AI
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.38"
$ws.Range("E31").Value = "  -3.41%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.72"
$ws.Range("E32").Value = "  -2.89%  "

# Row 33: PancakeSwap
$ws.Range("E33").Value = "  -2.58%  "

# Row 34: FirstDigitalUSD
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.03%  "

# Row 35: Kaspa
$ws.Range("E35").Value = "  -7.21%  "

# Row 36: ImmutableX
$ws.Range("E36").Value = "  -1.68%  "

# Row 37: NEARProtocol
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.84"
$ws.Range("E37").Value = "  -1.66%  "

# Row 38: Monero
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "154.21"
$ws.Range("E38").Value = "  -0.86%  "

# Row 39: EthereumClassic
$ws.Range("E39").Value = "  -2.70%  "

# Row 40: PolygonEcosystemToken
$ws.Range("E40").Value = "  -0.90%  "

# Row 41: RenderToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.23"
$ws.Range("E41").Value = "  -2.85%  "

# Row 42: Stacks
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.80"
$ws.Range("E42").Value = "  -2.94%  "

# Row 43: dogwifhat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  +0.25%  "

# Row 44: WhiteBITCoin
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.80"
$ws.Range("E44").Value = "  +1.67%  "

# Row 45: OKB
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.92"
$ws.Range("E45").Value = "  -1.07%  "

# Row 46: USDe
$ws.Range("E46").Value = "  -0.02%  "

# Row 47: Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.62"
$ws.Range("E47").Value = "  -1.25%  "

# Row 48: BabyDogeCoin
$ws.Range("E48").Value = "  +0.92%  "

# Row 49: Filecoin
$ws.Range("E49").Value = "  -0.34%  "

# Row 50: InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.61"
$ws.Range("E50").Value = "  +3.54%  "

# Row 51: Mantle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.616"
$ws.Range("E51").Value = "  -1.94%  "
